# Insert a new weekly record as row 63 in the "Arveja Verde" sheet.
# This pushes the previous rows 63-73 down to 64-74 (data unchanged),
# grows the used range from A1:R73 to A1:R74, and fills the newly
# inserted row 63 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63..73 down to 64..74, leaving a blank row 63 behind.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the latest weekly observation.
$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44505
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100112022
$ws.Cells.Item(63, 7).Value = "Arveja Verde"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 48
$ws.Cells.Item(63, 11).Value = 11000
$ws.Cells.Item(63, 12).Value = 13000
$ws.Cells.Item(63, 13).Value = 12042
$ws.Cells.Item(63, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 482
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
